$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected Table 1 domain analysis values (Frequency columns reformatted,
# Percent columns recomputed) per commit "Table 1 domain analysis corrected".
$ws.Range("B2").Value = 2629.0
$ws.Range("C2").Value = 0.12486114311188569
$ws.Range("E2").Value = 502.0
$ws.Range("F2").Value = 0.18507942308033346
$ws.Range("G2").Value = 1185.0
$ws.Range("H2").Value = 0.4898956146358516
$ws.Range("B3").Value = 3135.0
$ws.Range("C3").Value = 0.26255386707176337
$ws.Range("E3").Value = 319.0
$ws.Range("F3").Value = 0.0797996404560808
$ws.Range("G3").Value = 1419.0
$ws.Range("H3").Value = 0.495217650639752
$ws.Range("B4").Value = 3840.0
$ws.Range("C4").Value = 0.24490627733231402
$ws.Range("E4").Value = 672.0
$ws.Range("F4").Value = 0.15581981197122155
$ws.Range("G4").Value = 1830.0
$ws.Range("H4").Value = 0.5244569490000929
$ws.Range("B5").Value = 2059.0
$ws.Range("C5").Value = 0.06598153257734969
$ws.Range("E5").Value = 356.0
$ws.Range("F5").Value = 0.16169120180728133
$ws.Range("G5").Value = 841.0
$ws.Range("H5").Value = 0.43136326543532655
$ws.Range("B6").Value = 4382.0
$ws.Range("C6").Value = 0.30169717990668726
$ws.Range("E6").Value = 739.0
$ws.Range("F6").Value = 0.14891291937343276
$ws.Range("G6").Value = 2040.0
$ws.Range("H6").Value = 0.5126616125361327
$ws.Range("B7").Value = 8073.0
$ws.Range("C7").Value = 0.5117627120635093
$ws.Range("E7").Value = 1297.0
$ws.Range("F7").Value = 0.1354983682157096
$ws.Range("G7").Value = 3840.0
$ws.Range("H7").Value = 0.523752820413486
$ws.Range("B8").Value = 7999.0
$ws.Range("C8").Value = 0.4882372879364906
$ws.Range("E8").Value = 1296.0
$ws.Range("F8").Value = 0.140672320312189
$ws.Range("G8").Value = 3482.0
$ws.Range("H8").Value = 0.48016432037236295
$ws.Range("B9").Value = 2895.0
$ws.Range("C9").Value = 0.08140105460711192
$ws.Range("E9").Value = 377.0
$ws.Range("F9").Value = 0.12119582655744879
$ws.Range("G9").Value = 1166.0
$ws.Range("H9").Value = 0.41280301368474204
$ws.Range("B10").Value = 3295.0
$ws.Range("C10").Value = 0.11339982690459612
$ws.Range("E10").Value = 851.0
$ws.Range("F10").Value = 0.263793461420665
$ws.Range("G10").Value = 1527.0
$ws.Range("H10").Value = 0.4641894244684604
$ws.Range("B11").Value = 7779.0
$ws.Range("C11").Value = 0.6997905026376184
$ws.Range("E11").Value = 987.0
$ws.Range("F11").Value = 0.11557434011829944
$ws.Range("G11").Value = 3817.0
$ws.Range("H11").Value = 0.529225081865937
$ws.Range("B12").Value = 1371.0
$ws.Range("C12").Value = 0.04431119017169008
$ws.Range("E12").Value = 246.0
$ws.Range("F12").Value = 0.1624588145847037
$ws.Range("G12").Value = 519.0
$ws.Range("H12").Value = 0.4167798661601008
$ws.Range("B13").Value = 732.0
$ws.Range("C13").Value = 0.06109742567898335
$ws.Range("E13").Value = 132.0
$ws.Range("F13").Value = 0.16659755285831082
$ws.Range("G13").Value = 293.0
$ws.Range("H13").Value = 0.448551759142935
$ws.Range("B14").Value = 1004.0
$ws.Range("C14").Value = 0.05747827000119642
$ws.Range("E14").Value = 192.0
$ws.Range("F14").Value = 0.16712454268678914
$ws.Range("G14").Value = 389.0
$ws.Range("H14").Value = 0.4261183170491705
$ws.Range("B15").Value = 7999.0
$ws.Range("C15").Value = 0.4882372879364907
$ws.Range("E15").Value = 1296.0
$ws.Range("F15").Value = 0.140672320312189
$ws.Range("G15").Value = 3482.0
$ws.Range("H15").Value = 0.48016432037236295
$ws.Range("B16").Value = 6545.0
$ws.Range("C16").Value = 0.4110909055852731
$ws.Range("E16").Value = 1037.0
$ws.Range("F16").Value = 0.13419117579856057
$ws.Range("G16").Value = 3188.0
$ws.Range("H16").Value = 0.5306338017094295
$ws.Range("B17").Value = 524.0
$ws.Range("C17").Value = 0.04319353647703991
$ws.Range("E17").Value = 68.0
$ws.Range("F17").Value = 0.10586304778593553
$ws.Range("G17").Value = 263.0
$ws.Range("H17").Value = 0.5880940858461586
$ws.Range("B18").Value = 4191.0
$ws.Range("C18").Value = 0.2668761193923869
$ws.Range("E18").Value = 840.0
$ws.Range("F18").Value = 0.1867148477516125
$ws.Range("G18").Value = 2022.0
$ws.Range("H18").Value = 0.5299426217271924
$ws.Range("B19").Value = 11881.0
$ws.Range("C19").Value = 0.7331238806076127
$ws.Range("E19").Value = 1753.0
$ws.Range("F19").Value = 0.12031952446696824
$ws.Range("G19").Value = 5300.0
$ws.Range("H19").Value = 0.4924662840631326
$ws.Range("B20").Value = 1191.0
$ws.Range("C20").Value = 0.07141427156719986
$ws.Range("E20").Value = 212.0
$ws.Range("F20").Value = 0.15811668149210137
$ws.Range("G20").Value = 486.0
$ws.Range("H20").Value = 0.4535537831821629
$ws.Range("B21").Value = 7999.0
$ws.Range("C21").Value = 0.4882372879364907
$ws.Range("E21").Value = 1296.0
$ws.Range("F21").Value = 0.140672320312189
$ws.Range("G21").Value = 3482.0
$ws.Range("H21").Value = 0.48016432037236295
$ws.Range("B22").Value = 6616.0
$ws.Range("C22").Value = 0.4196533565174644
$ws.Range("E22").Value = 1053.0
$ws.Range("F22").Value = 0.1337839473476835
$ws.Range("G22").Value = 3190.0
$ws.Range("H22").Value = 0.5292337080994438
$ws.Range("B23").Value = 266.0
$ws.Range("C23").Value = 0.02069508397884525
$ws.Range("E23").Value = 32.0
$ws.Range("F23").Value = 0.0922032459515263
$ws.Range("G23").Value = 164.0
$ws.Range("H23").Value = 0.6547178499619064
$ws.Range("B24").Value = 10262.0
$ws.Range("C24").Value = 0.6587112722226114
$ws.Range("E24").Value = 1499.0
$ws.Range("F24").Value = 0.1261754722221648
$ws.Range("G24").Value = 4702.0
$ws.Range("H24").Value = 0.5072758279664283
$ws.Range("B25").Value = 5810.0
$ws.Range("C25").Value = 0.3412887277773887
$ws.Range("E25").Value = 1094.0
$ws.Range("F25").Value = 0.1608998184941224
$ws.Range("G25").Value = 2620.0
$ws.Range("H25").Value = 0.49319278243276726
$ws.Range("B26").Value = 4442.0
$ws.Range("C26").Value = 0.3236744170748077
$ws.Range("E26").Value = 695.0
$ws.Range("F26").Value = 0.13479304388939517
$ws.Range("G26").Value = 3270.0
$ws.Range("H26").Value = 0.7692280837800479
$ws.Range("B27").Value = 3409.0
$ws.Range("C27").Value = 0.15919464430513572
$ws.Range("E27").Value = 594.0
$ws.Range("F27").Value = 0.15434045425639598
$ws.Range("G27").Value = 1188.0
$ws.Range("H27").Value = 0.35881286903998816
$ws.Range("B28").Value = 8221.0
$ws.Range("C28").Value = 0.5171309386200567
$ws.Range("E28").Value = 1304.0
$ws.Range("F28").Value = 0.13502321735951536
$ws.Range("G28").Value = 2864.0
$ws.Range("H28").Value = 0.3796871203729662
$ws.Range("B29").Value = 4316.0
$ws.Range("C29").Value = 0.1897291484489276
$ws.Range("E29").Value = 839.0
$ws.Range("F29").Value = 0.1872257421535215
$ws.Range("G29").Value = 1778.0
$ws.Range("H29").Value = 0.45477657312446507
$ws.Range("B30").Value = 3962.0
$ws.Range("C30").Value = 0.18997926591229702
$ws.Range("E30").Value = 683.0
$ws.Range("F30").Value = 0.1691324883782244
$ws.Range("G30").Value = 1794.0
$ws.Range("H30").Value = 0.5074442518790123
$ws.Range("B31").Value = 7794.0
$ws.Range("C31").Value = 0.6202915856387755
$ws.Range("E31").Value = 1071.0
$ws.Range("F31").Value = 0.1134919334581784
$ws.Range("G31").Value = 3750.0
$ws.Range("H31").Value = 0.5155330137327362
$ws.Range("B32").Value = 8741.0
$ws.Range("C32").Value = 0.4975309865055008
$ws.Range("E32").Value = 1143.0
$ws.Range("F32").Value = 0.10819496651045324
$ws.Range("B33").Value = 7322.0
$ws.Range("C33").Value = 0.5024690134944994
$ws.Range("E33").Value = 1448.0
$ws.Range("F33").Value = 0.1674067851789359
$ws.Range("B34").Value = 13460.0
$ws.Range("C34").Value = 0.8619747054191564
$ws.Range("G34").Value = 5862.0
$ws.Range("H34").Value = 0.48516092861343796
$ws.Range("B35").Value = 2593.0
$ws.Range("C35").Value = 0.13802529458084364
$ws.Range("G35").Value = 1448.0
$ws.Range("H35").Value = 0.6096450559820747
$ws.Range("B36").Value = 612.0
$ws.Range("C36").Value = 0.04055348062536871
$ws.Range("E36").Value = 153.0
$ws.Range("F36").Value = 0.21491655551873742
$ws.Range("G36").Value = 548.0
$ws.Range("H36").Value = 0.9046874261574475
$ws.Range("B37").Value = 14659.0
$ws.Range("C37").Value = 0.9020236647277617
$ws.Range("E37").Value = 2246.0
$ws.Range("F37").Value = 0.12978994006035618
$ws.Range("G37").Value = 6013.0
$ws.Range("H37").Value = 0.4552273107625041
$ws.Range("B38").Value = 801.0
$ws.Range("C38").Value = 0.057422854646869116
$ws.Range("E38").Value = 194.0
$ws.Range("F38").Value = 0.21314409695001446
$ws.Range("G38").Value = 761.0
$ws.Range("H38").Value = 0.960469825911913
